$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.396.59"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").Value = "1.670.81"
$ws.Range("E3").Value = "  +3.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5296"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.68%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.670.95"
$ws.Range("E12").Value = "  +3.95%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.485"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5548"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.02%  "

$ws.Range("D15").Value = "0.0₅8332"
$ws.Range("E15").Value = "  +7.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.89%  "

$ws.Range("D17").Value = "26.476.16"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.769"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.316"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.415"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.61%  "

$ws.Range("E28").Value = "  +5.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.272"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.614"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.439"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.685"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.004"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.429"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.770"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5697"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01636"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.75%  "

$ws.Range("D40").Value = "1.066.38"
$ws.Range("E40").Value = "  +4.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8598"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.30%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("D44").Value = "1.827.66"
$ws.Range("E44").Value = "  +3.81%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈108"
$ws.Range("E45").Value = "  +2.58%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.136"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4242"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.010"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.57%  "

